$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a "Benchmark" header label in C1 (previously empty) for the new
# test-workflow column used by the arena evaluation sheet.
$ws.Range("C1").Value = "Benchmark"
